$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The second data row (row 2) was a duplicate of what is now needed one row
# later: shift the block of data rows 2-5 down into rows 3-6 (overwriting
# the old row 6), then clear out the now-vacated row 2 so it disappears
# from the sheet entirely (matching the diff, which drops the <row r="2">
# element and re-numbers the remaining data down into rows 3-6).
$ws.Range("A2:F5").Copy($ws.Range("A3:F6"))
$ws.Range("A2:F2").ClearContents()

Write-Output "Shifted rows 2:5 down into 3:6 and cleared row 2"
